$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 12.1164
$ws.Range("A8").Value = -21.09
$ws.Range("A10").Value = -20.50069999999997
$ws.Range("A12").Value = -22.38400000000004
$ws.Range("B13").Value = 6.560599999999996
$ws.Range("A18").Value = -22.41320000000004
$ws.Range("E20").Value = 12.3767
$ws.Range("A25").Value = -22.23540000000003
